# Generate Report for Handoff
# Update Priority from "low" to "ht" and refresh handoff timestamps
# for the four file rows (4-7) that had not yet been handed off.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsOverview = $wb.Worksheets.Item("Overview")

# zh-cn sheet: Priority (E) -> "ht", Latest Handoff Datetime (H) -> new timestamp
foreach ($row in 4..7) {
    $wsZhCn.Cells.Item($row, 5).Value = "ht"
    $wsZhCn.Cells.Item($row, 8).Value = "2016-08-20 18:43:14"
}

# de-de sheet: Priority (E) -> "ht", Latest Handoff Datetime (H) -> new timestamp
foreach ($row in 4..7) {
    $wsDeDe.Cells.Item($row, 5).Value = "ht"
    $wsDeDe.Cells.Item($row, 8).Value = "2016-08-20 18:43:18"
}

# Overview sheet: Latest HO Xliff Generate Date (G) -> new timestamp
foreach ($row in 4..7) {
    $wsOverview.Cells.Item($row, 7).Value = "2016-08-20 18:43:18"
}
